# Update parametric survival model parameter tables on the first three
# worksheets (PFS caba, PFS mito, OS caba). The underlying R data structure
# used to build these tables had changed shape, so the function was pulling
# the wrong values; this corrects the AIC, intercept, and log(scale) rows.

$wb = $excel.ActiveWorkbook

function Set-SheetValues($Sheet, $Row2, $Row3, $Row4) {
    # Columns B, C, D, E correspond to exponential, weibull, lognormal, loglogistic
    $cols = @("B", "C", "D", "E")

    for ($i = 0; $i -lt 4; $i++) {
        $Sheet.Range("$($cols[$i])2").Value = $Row2[$i]
        $Sheet.Range("$($cols[$i])3").Value = $Row3[$i]
    }

    # Row 4 (log(scale)) only has values for columns C, D, E
    $Sheet.Range("C4").Value = $Row4[0]
    $Sheet.Range("D4").Value = $Row4[1]
    $Sheet.Range("E4").Value = $Row4[2]
}

# Sheet 1: PFS caba
$ws1 = $wb.Worksheets.Item("PFS caba")
Set-SheetValues $ws1 `
    @(1885.1071, 1885.5943, 1886.6016, 1895.3382) `
    @(1.4561, 1.4383, 0.9282, 0.9382) `
    @(0.0566, 0.1527, -0.3479)

# Sheet 2: PFS mito
$ws2 = $wb.Worksheets.Item("PFS mito")
Set-SheetValues $ws2 `
    @(1754.2742, 1748.7375, 1691.7935, 1695.7627) `
    @(1.1368, 1.0851, 0.584, 0.5222) `
    @(0.1132, 0.0635, -0.4662)

# Sheet 3: OS caba
$ws3 = $wb.Worksheets.Item("OS caba")
Set-SheetValues $ws3 `
    @(1987.5033, 1939.5375, 1952.9809, 1938.9274) `
    @(3.0551, 3.0125, 2.705, 2.723) `
    @(-0.405, -0.1, -0.68)

# Sheet 4 (OS mito) is intentionally left unchanged.
